# This workbook's rows 5-30 (the "artfynd" observation rows) were
# re-ordered upstream: each destination row now carries the data that
# used to live in a different source row (columns A..AY). Build that
# row -> source-row mapping, then permute the block in one shot via a
# Value2 array round-trip so every column (ids, coordinates, counts,
# dates, ...) moves together and precision is preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (both absolute sheet row numbers)
$rowMap = @{
    5  = 25
    6  = 5
    7  = 6
    8  = 26
    9  = 7
    10 = 8
    11 = 9
    12 = 10
    13 = 11
    14 = 12
    15 = 27
    16 = 28
    17 = 13
    18 = 29
    19 = 14
    20 = 15
    21 = 16
    22 = 17
    23 = 18
    24 = 30
    25 = 19
    26 = 20
    27 = 21
    28 = 22
    29 = 23
    30 = 24
}

$firstRow = 5
$lastRow = 30
$firstCol = 1   # A
$lastCol = 51   # AY

# Snapshot the original block before any writes (source rows must be
# read from the pristine data, not from cells already overwritten).
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$orig = $srcRange.Value2

$numRows = $lastRow - $firstRow + 1
$numCols = $lastCol - $firstCol + 1

# Build the permuted block in memory. New-Object '[,]' arrays are
# 0-based, unlike the 1-based COM Value2 array read above.
$dest = New-Object 'object[,]' $numRows, $numCols
foreach ($destRow in $firstRow..$lastRow) {
    $srcRow = $rowMap[$destRow]
    $destIdx = $destRow - $firstRow
    $srcIdx = $srcRow - $firstRow + 1
    for ($c = 1; $c -le $numCols; $c++) {
        $dest[$destIdx, ($c - 1)] = $orig[$srcIdx, $c]
    }
}

# Some columns hold text that Value2 would otherwise coerce to a
# number/date on write-back:
#  - I ("Antal") stores digit strings like "8" as text, not numbers.
#  - Y/AA hold literal "YYYY-MM-DD" text that would auto-convert to a
#    real date serial.
# Every row 5-30 already has a (possibly empty) cell in these columns,
# so forcing text format here cannot spuriously materialise a cell
# that should stay fully absent (unlike e.g. K/L/N/AF).
foreach ($c in @(9, 25, 27)) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $c), $ws.Cells.Item($lastRow, $c))
    $colRange.NumberFormat = "@"
}

$destRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$destRange.Value2 = $dest
